# Auto-generated edit script applying the Omega_Profits financial-value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (102 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 270
$ws.Range("I4").Value = 140
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 140
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -26
$ws.Range("N4").Value = -628
$ws.Range("H15").Value = 399.13333
$ws.Range("I15").Value = 399.13333
$ws.Range("K15").Value = 1197.39999
$ws.Range("M15").Value = -1028.39999
$ws.Range("H38").Value = 3433.4
$ws.Range("I38").Value = 160.1
$ws.Range("K38").Value = 480.3
$ws.Range("M38").Value = -108.3
$ws.Range("H43").Value = 4269.8335
$ws.Range("J43").Value = 4149
$ws.Range("L43").Value = 4149
$ws.Range("N43").Value = -4287
$ws.Range("H64").Value = 7455.1665
$ws.Range("I64").Value = 5582
$ws.Range("K64").Value = 5582
$ws.Range("M64").Value = -5334
$ws.Range("H67").Value = 7455.1665
$ws.Range("I67").Value = 5582
$ws.Range("K67").Value = 5582
$ws.Range("M67").Value = -4724
$ws.Range("H69").Value = 16625
$ws.Range("J69").Value = 17642.857
$ws.Range("L69").Value = 52928.571
$ws.Range("N69").Value = -54676.571
$ws.Range("H72").Value = 16625
$ws.Range("J72").Value = 17642.857
$ws.Range("L72").Value = 158785.713
$ws.Range("N72").Value = -167521.713
$ws.Range("H92").Value = 1311.5555
$ws.Range("I92").Value = 1340.2307
$ws.Range("K92").Value = 1340.2307
$ws.Range("M92").Value = -92.23070000000007
$ws.Range("H93").Value = 52500
$ws.Range("J93").Value = 52500
$ws.Range("L93").Value = 52500
$ws.Range("N93").Value = -57492
$ws.Range("H96").Value = 1272.7
$ws.Range("I96").Value = 783
$ws.Range("J96").Value = 2007.25
$ws.Range("K96").Value = 2349
$ws.Range("L96").Value = 6021.75
$ws.Range("M96").Value = -976
$ws.Range("N96").Value = -8767.75
$ws.Range("H99").Value = 1308.4546
$ws.Range("I99").Value = 1060.8572
$ws.Range("J99").Value = 1741.75
$ws.Range("K99").Value = 3182.5716
$ws.Range("L99").Value = 5225.25
$ws.Range("M99").Value = -1684.5716
$ws.Range("N99").Value = -8221.25
$ws.Range("H101").Value = 503.83334
$ws.Range("J101").Value = 500.66666
$ws.Range("L101").Value = 1501.99998
$ws.Range("N101").Value = -4745.999980000001
$ws.Range("H113").Value = 2655.1667
$ws.Range("I113").Value = 2655.1667
$ws.Range("K113").Value = 2655.1667
$ws.Range("M113").Value = 598.8332999999998
$ws.Range("H116").Value = 4874.25
$ws.Range("I116").Value = 4774.75
$ws.Range("K116").Value = 4774.75
$ws.Range("M116").Value = -1332.75
$ws.Range("H121").Value = 4220.8335
$ws.Range("J121").Value = 4220.8335
$ws.Range("L121").Value = 12662.5005
$ws.Range("N121").Value = -16156.5005
$ws.Range("H125").Value = 1032
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null
$ws.Range("H132").Value = 3121.8
$ws.Range("I132").Value = 3206.2222
$ws.Range("J132").Value = 2362
$ws.Range("K132").Value = 9618.6666
$ws.Range("L132").Value = 7086
$ws.Range("M132").Value = -7088.6666
$ws.Range("N132").Value = -12146
$ws.Range("H137").Value = 2818.6667
$ws.Range("J137").Value = 3116
$ws.Range("L137").Value = 9348
$ws.Range("N137").Value = -14448
$ws.Range("H138").Value = 3136.1562
$ws.Range("I138").Value = 1609.8462
$ws.Range("J138").Value = 4180.4736
$ws.Range("K138").Value = 4829.5386
$ws.Range("L138").Value = 12541.4208
$ws.Range("M138").Value = 310.4614000000001
$ws.Range("N138").Value = -22821.4208
$ws.Range("H141").Value = 6428.4287
$ws.Range("I141").Value = 6125
$ws.Range("J141").Value = 6833
$ws.Range("K141").Value = 18375
$ws.Range("L141").Value = 20499
$ws.Range("M141").Value = -13195
$ws.Range("N141").Value = -30859

# ---- Sheet: ARM (46 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3537.6775
$ws.Range("I2").Value = 3495.4644
$ws.Range("K2").Value = 3495.4644
$ws.Range("M2").Value = -3382.4644
$ws.Range("H32").Value = 9461.596
$ws.Range("I32").Value = 830.0909
$ws.Range("K32").Value = 830.0909
$ws.Range("M32").Value = -543.0909
$ws.Range("H45").Value = 11102
$ws.Range("I45").Value = 16279.6
$ws.Range("J45").Value = 4630
$ws.Range("K45").Value = 16279.6
$ws.Range("L45").Value = 4630
$ws.Range("M45").Value = -15902.6
$ws.Range("N45").Value = -5384
$ws.Range("H87").Value = 25000
$ws.Range("I87").Value = 25000
$ws.Range("K87").Value = 25000
$ws.Range("M87").Value = -23752
$ws.Range("H90").Value = 25000
$ws.Range("I90").Value = 25000
$ws.Range("K90").Value = 75000
$ws.Range("M90").Value = -68760
$ws.Range("H106").Value = 41663.332
$ws.Range("J106").Value = 41663.332
$ws.Range("L106").Value = 41663.332
$ws.Range("N106").Value = -44187.332
$ws.Range("H110").Value = 1258.0834
$ws.Range("I110").Value = 1190.6818
$ws.Range("K110").Value = 1190.6818
$ws.Range("M110").Value = 854.3181999999999
$ws.Range("H111").Value = 42660.25
$ws.Range("J111").Value = 42660.25
$ws.Range("L111").Value = 42660.25
$ws.Range("N111").Value = -50840.25
$ws.Range("H116").Value = 3537.6775
$ws.Range("I116").Value = 3495.4644
$ws.Range("K116").Value = 3495.4644
$ws.Range("M116").Value = -1201.4644
$ws.Range("H132").Value = 3439.4783
$ws.Range("I132").Value = 2919.5264
$ws.Range("J132").Value = 5909.25
$ws.Range("K132").Value = 8758.5792
$ws.Range("L132").Value = 17727.75
$ws.Range("M132").Value = -6228.5792
$ws.Range("N132").Value = -22787.75

# ---- Sheet: BSM (40 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3537.6775
$ws.Range("I3").Value = 3495.4644
$ws.Range("K3").Value = 3495.4644
$ws.Range("M3").Value = -3381.4644
$ws.Range("H20").Value = 1544.125
$ws.Range("I20").Value = 1762.3529
$ws.Range("K20").Value = 1762.3529
$ws.Range("M20").Value = -1515.3529
$ws.Range("H80").Value = 1144.8695
$ws.Range("I80").Value = 1383
$ws.Range("J80").Value = 885.0909
$ws.Range("K80").Value = 1383
$ws.Range("L80").Value = 885.0909
$ws.Range("M80").Value = -385
$ws.Range("N80").Value = -2881.0909
$ws.Range("H83").Value = 1144.8695
$ws.Range("I83").Value = 1383
$ws.Range("J83").Value = 885.0909
$ws.Range("K83").Value = 6915
$ws.Range("L83").Value = 4425.4545
$ws.Range("M83").Value = -1923
$ws.Range("N83").Value = -14409.4545
$ws.Range("H105").Value = 2632.7778
$ws.Range("I105").Value = 2142.2856
$ws.Range("K105").Value = 2142.2856
$ws.Range("M105").Value = -395.2856000000002
$ws.Range("H107").Value = 787.8
$ws.Range("I107").Value = 734.875
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 734.875
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 1185.125
$ws.Range("N107").Value = -4839.5
$ws.Range("H134").Value = 2412.4375
$ws.Range("I134").Value = 2184.6428
$ws.Range("J134").Value = 4007
$ws.Range("K134").Value = 6553.928400000001
$ws.Range("L134").Value = 12021
$ws.Range("M134").Value = -4018.928400000001
$ws.Range("N134").Value = -17091

# ---- Sheet: CRP (82 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 36961.79
$ws.Range("I16").Value = 18140.84
$ws.Range("J16").Value = 73155.92
$ws.Range("K16").Value = 18140.84
$ws.Range("L16").Value = 73155.92
$ws.Range("M16").Value = -17853.84
$ws.Range("N16").Value = -73729.92
$ws.Range("H31").Value = 6031.1055
$ws.Range("I31").Value = 6956.4
$ws.Range("K31").Value = 6956.4
$ws.Range("M31").Value = -6661.4
$ws.Range("H34").Value = 6031.1055
$ws.Range("I34").Value = 6956.4
$ws.Range("K34").Value = 6956.4
$ws.Range("M34").Value = -6754.4
$ws.Range("H50").Value = 22000
$ws.Range("J50").Value = 22000
$ws.Range("L50").Value = 22000
$ws.Range("N50").Value = -23250
$ws.Range("H58").Value = 5599.8335
$ws.Range("I58").Value = 5811.769
$ws.Range("K58").Value = 5811.769
$ws.Range("M58").Value = -5608.769
$ws.Range("H62").Value = 15000
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 15000
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81240
$ws.Range("H95").Value = 28290.666
$ws.Range("J95").Value = 28290.666
$ws.Range("L95").Value = 28290.666
$ws.Range("N95").Value = -33782.666
$ws.Range("H96").Value = 28037.223
$ws.Range("J96").Value = 28037.223
$ws.Range("L96").Value = 28037.223
$ws.Range("N96").Value = -33529.223
$ws.Range("H99").Value = 2928323.5
$ws.Range("I99").Value = 3972661.5
$ws.Range("J99").Value = 4177.4
$ws.Range("K99").Value = 3972661.5
$ws.Range("L99").Value = 4177.4
$ws.Range("M99").Value = -3971163.5
$ws.Range("N99").Value = -7173.4
$ws.Range("H113").Value = 36961.79
$ws.Range("I113").Value = 18140.84
$ws.Range("J113").Value = 73155.92
$ws.Range("K113").Value = 18140.84
$ws.Range("L113").Value = 73155.92
$ws.Range("M113").Value = -15970.84
$ws.Range("N113").Value = -77495.92
$ws.Range("H122").Value = 3541.2144
$ws.Range("I122").Value = 3582.8462
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10748.5386
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8298.5386
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 2928323.5
$ws.Range("I126").Value = 3972661.5
$ws.Range("J126").Value = 4177.4
$ws.Range("K126").Value = 11917984.5
$ws.Range("L126").Value = 12532.2
$ws.Range("M126").Value = -11915514.5
$ws.Range("N126").Value = -17472.2
$ws.Range("H132").Value = 6242.5
$ws.Range("I132").Value = 4854.647
$ws.Range("J132").Value = 8864
$ws.Range("K132").Value = 14563.941
$ws.Range("L132").Value = 26592
$ws.Range("M132").Value = -12033.941
$ws.Range("N132").Value = -31652
$ws.Range("H134").Value = 3565
$ws.Range("I134").Value = 3602.6667
$ws.Range("K134").Value = 10808.0001
$ws.Range("M134").Value = -8273.000100000001
$ws.Range("H136").Value = 5599.8335
$ws.Range("I136").Value = 5811.769
$ws.Range("K136").Value = 17435.307
$ws.Range("M136").Value = -14885.307

# ---- Sheet: CUL (32 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 606.5484
$ws.Range("I2").Value = 908.05884
$ws.Range("K2").Value = 5448.35304
$ws.Range("M2").Value = -5335.35304
$ws.Range("H5").Value = 1619.7142
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 15000
$ws.Range("N5").Value = -15224
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H26").Value = 90
$ws.Range("J26").Value = 90
$ws.Range("L26").Value = 270
$ws.Range("N26").Value = -846
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H131").Value = 2008.9524
$ws.Range("I131").Value = 1136.875
$ws.Range("K131").Value = 3410.625
$ws.Range("M131").Value = 1629.375
$ws.Range("H135").Value = 1619.7142
$ws.Range("J135").Value = 5000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -50070
$ws.Range("H139").Value = 6008.615
$ws.Range("I139").Value = 4018.1428
$ws.Range("K139").Value = 12054.4284
$ws.Range("M139").Value = -6914.428400000001

# ---- Sheet: GSM (57 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 221.42857
$ws.Range("I2").Value = 242.5
$ws.Range("J2").Value = 95
$ws.Range("K2").Value = 242.5
$ws.Range("L2").Value = 95
$ws.Range("M2").Value = -129.5
$ws.Range("N2").Value = -321
$ws.Range("H70").Value = 3801.9167
$ws.Range("I70").Value = 3649.9412
$ws.Range("J70").Value = 4171
$ws.Range("K70").Value = 3649.9412
$ws.Range("L70").Value = 4171
$ws.Range("M70").Value = -3379.9412
$ws.Range("N70").Value = -4711
$ws.Range("H73").Value = 3801.9167
$ws.Range("I73").Value = 3649.9412
$ws.Range("J73").Value = 4171
$ws.Range("K73").Value = 3649.9412
$ws.Range("L73").Value = 4171
$ws.Range("M73").Value = -2713.9412
$ws.Range("N73").Value = -6043
$ws.Range("H80").Value = 7257.074
$ws.Range("I80").Value = 7105
$ws.Range("J80").Value = 7378.7334
$ws.Range("K80").Value = 7105
$ws.Range("L80").Value = 7378.7334
$ws.Range("M80").Value = -6107
$ws.Range("N80").Value = -9374.733400000001
$ws.Range("H83").Value = 7257.074
$ws.Range("I83").Value = 7105
$ws.Range("J83").Value = 7378.7334
$ws.Range("K83").Value = 35525
$ws.Range("L83").Value = 36893.667
$ws.Range("M83").Value = -30533
$ws.Range("N83").Value = -46877.667
$ws.Range("H92").Value = 8498.75
$ws.Range("J92").Value = 8498.75
$ws.Range("L92").Value = 8498.75
$ws.Range("N92").Value = -12242.75
$ws.Range("H102").Value = 2731.739
$ws.Range("I102").Value = 2264.7896
$ws.Range("J102").Value = 4949.75
$ws.Range("K102").Value = 2264.7896
$ws.Range("L102").Value = 4949.75
$ws.Range("M102").Value = -642.7896000000001
$ws.Range("N102").Value = -8193.75
$ws.Range("H126").Value = 3943.8125
$ws.Range("I126").Value = 2454.7273
$ws.Range("J126").Value = 7219.8
$ws.Range("K126").Value = 7364.1819
$ws.Range("L126").Value = 21659.4
$ws.Range("M126").Value = -4894.1819
$ws.Range("N126").Value = -26599.4
$ws.Range("H132").Value = 3695.8
$ws.Range("I132").Value = 3489.2942
$ws.Range("K132").Value = 10467.8826
$ws.Range("M132").Value = -7937.882599999999

# ---- Sheet: LTW (61 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H16").Value = 2766
$ws.Range("I16").Value = 1630
$ws.Range("J16").Value = 5416.6665
$ws.Range("K16").Value = 1630
$ws.Range("L16").Value = 5416.6665
$ws.Range("M16").Value = -1460
$ws.Range("N16").Value = -5756.6665
$ws.Range("H22").Value = 2870.3333
$ws.Range("I22").Value = 2305.5
$ws.Range("K22").Value = 2305.5
$ws.Range("M22").Value = -2010.5
$ws.Range("H23").Value = 4000000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
$ws.Range("H27").Value = 2870.3333
$ws.Range("I27").Value = 2305.5
$ws.Range("K27").Value = 2305.5
$ws.Range("M27").Value = -2198.5
$ws.Range("H40").Value = 10710.52
$ws.Range("I40").Value = 7605.7646
$ws.Range("J40").Value = 17308.125
$ws.Range("K40").Value = 7605.7646
$ws.Range("L40").Value = 17308.125
$ws.Range("M40").Value = -7469.7646
$ws.Range("N40").Value = -17580.125
$ws.Range("H55").Value = 771.125
$ws.Range("I55").Value = 466.1
$ws.Range("J55").Value = 1279.5
$ws.Range("K55").Value = 466.1
$ws.Range("L55").Value = 1279.5
$ws.Range("M55").Value = -293.1
$ws.Range("N55").Value = -1625.5
$ws.Range("H61").Value = 2955.5
$ws.Range("J61").Value = 2763.75
$ws.Range("L61").Value = 2763.75
$ws.Range("N61").Value = -3167.75
$ws.Range("H103").Value = 62500
$ws.Range("J103").Value = 62500
$ws.Range("L103").Value = 62500
$ws.Range("N103").Value = -64844
$ws.Range("H108").Value = 45763.668
$ws.Range("J108").Value = 45763.668
$ws.Range("L108").Value = 45763.668
$ws.Range("N108").Value = -53443.668
$ws.Range("H113").Value = 2955.5
$ws.Range("J113").Value = 2763.75
$ws.Range("L113").Value = 2763.75
$ws.Range("N113").Value = -7103.75
$ws.Range("H122").Value = 5699.3076
$ws.Range("I122").Value = 4454.6665
$ws.Range("K122").Value = 13363.9995
$ws.Range("M122").Value = -10913.9995

# ---- Sheet: WVR (40 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H96").Value = 2502.6667
$ws.Range("I96").Value = 2502.6667
$ws.Range("K96").Value = 2502.6667
$ws.Range("M96").Value = -1129.6667
$ws.Range("H107").Value = 2144.2917
$ws.Range("I107").Value = 1731
$ws.Range("J107").Value = 2439.5
$ws.Range("K107").Value = 5193
$ws.Range("L107").Value = 7318.5
$ws.Range("M107").Value = -3273
$ws.Range("N107").Value = -11158.5
$ws.Range("H122").Value = 4492
$ws.Range("I122").Value = 3467.182
$ws.Range("K122").Value = 10401.546
$ws.Range("M122").Value = -7951.545999999998
$ws.Range("H132").Value = 5883.5625
$ws.Range("I132").Value = 5883.5625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17650.6875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15120.6875
$ws.Range("N132").Value = $null
$ws.Range("H136").Value = 12392.8125
$ws.Range("I136").Value = 11990.417
$ws.Range("J136").Value = 13600
$ws.Range("K136").Value = 35971.251
$ws.Range("L136").Value = 40800
$ws.Range("M136").Value = -33421.251
$ws.Range("N136").Value = -45900
$ws.Range("H138").Value = 143997.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 143997.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 143997.5
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = -154277.5
